# Fill in the missing "Release 4.0" progress figures on row 14 and mark the
# task as accepted ("Oui"), matching what was already done on row 13.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H14").Value = 4.03
$ws.Range("I14").Value = 4.03
$ws.Range("J14").Value = "Oui"

# Move the active selection to the cell that was just finished editing.
$ws.Range("J14").Select()
